# Update countries & provincias Spain
# Refresh the COVID country stats table on sheet "Pais" with the
# 23:08 data pull (was 21:51), which re-sorts a few adjacent rows
# (descending by "Casos totales") so some country names swap rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "datos actualizados" timestamp
$ws.Range("A1").Value = 'Datos actualizados a 15 de Julio de 2020 a las 23:08'

# Row 4: 'Estados Unidos' -> 'Estados Unidos'
$ws.Range("B4").Value = 3602202
$ws.Range("C4").Value = 57125
$ws.Range("D4").Value = 1626395
$ws.Range("E4").Value = 1835978
$ws.Range("G4").Value = 686
$ws.Range("H4").Value = 139829

# Row 8: 'Peru' -> 'Peru'
$ws.Range("B8").Value = 337724
$ws.Range("C8").Value = 3857
$ws.Range("D8").Value = 226400
$ws.Range("E8").Value = 98907
$ws.Range("G8").Value = 188
$ws.Range("H8").Value = 12417

# Row 19: 'Alemania' -> 'Alemania'
$ws.Range("B19").Value = 201252
$ws.Range("C19").Value = 486
$ws.Range("E19").Value = 6104
$ws.Range("G19").Value = 4
$ws.Range("H19").Value = 9148

# Row 46: 'Israel' -> 'Israel'
$ws.Range("B46").Value = 44188
$ws.Range("C46").Value = 1828
$ws.Range("D46").Value = 19989
$ws.Range("E46").Value = 23823
$ws.Range("G46").Value = 5
$ws.Range("H46").Value = 376

# Row 50: 'Barein' -> 'Barein'
$ws.Range("E50").Value = 4208
$ws.Range("G50").Value = 6
$ws.Range("H50").Value = 117

# Row 70: 'Dinamarca' -> 'Costa de Marfil'
$ws.Range("A70").Value = 'Costa de Marfil'
$ws.Range("B70").Value = 13403
$ws.Range("C70").Value = 366
$ws.Range("D70").Value = 7146
$ws.Range("E70").Value = 6170
$ws.Range("H70").Value = 87

# Row 71: 'Costa de Marfil' -> 'Dinamarca'
$ws.Range("A71").Value = 'Dinamarca'
$ws.Range("B71").Value = 13092
$ws.Range("C71").Value = 31
$ws.Range("D71").Value = 12182
$ws.Range("E71").Value = 300
$ws.Range("H71").Value = 610

# Row 97: 'Republica de Yibuti' -> 'Republica de Yibuti'
$ws.Range("B97").Value = 4985
$ws.Range("C97").Value = 6
$ws.Range("D97").Value = 4765
$ws.Range("E97").Value = 164

# Row 130: 'Suazilandia' -> 'Ruanda'
$ws.Range("A130").Value = 'Ruanda'
$ws.Range("B130").Value = 1435
$ws.Range("C130").Value = 19
$ws.Range("D130").Value = 752
$ws.Range("E130").Value = 679
$ws.Range("H130").Value = 4

# Row 131: 'Ruanda' -> 'Suazilandia'
$ws.Range("A131").Value = 'Suazilandia'
$ws.Range("B131").Value = 1434
$ws.Range("D131").Value = 695
$ws.Range("E131").Value = 719
$ws.Range("H131").Value = 20

# Row 139: 'Zimbabue' -> 'Zimbabue'
$ws.Range("B139").Value = 1089
$ws.Range("C139").Value = 25
$ws.Range("D139").Value = 395
$ws.Range("E139").Value = 674

# Row 149: 'Surinam' -> 'Surinam'
$ws.Range("B149").Value = 834
$ws.Range("C149").Value = 33
$ws.Range("D149").Value = 573
$ws.Range("E149").Value = 243

# Row 151: 'Santo Tome y Principe' -> 'Togo'
$ws.Range("A151").Value = 'Togo'
$ws.Range("B151").Value = 740
$ws.Range("C151").Value = 9
$ws.Range("D151").Value = 534
$ws.Range("E151").Value = 191
$ws.Range("H151").Value = 15

# Row 152: 'Togo' -> 'Santo Tome y Principe'
$ws.Range("A152").Value = 'Santo Tome y Principe'
$ws.Range("B152").Value = 732
$ws.Range("D152").Value = 286
$ws.Range("E152").Value = 432
$ws.Range("H152").Value = 14

# Row 209: 'Groenlandia' -> 'Islas Malvinas'
$ws.Range("A209").Value = 'Islas Malvinas'

# Row 210: 'Islas Malvinas' -> 'Groenlandia'
$ws.Range("A210").Value = 'Groenlandia'

